$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "CasesTab" row label to "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Update the selected cell from B3 to A2
$ws.Range("A2").Select()
